# Atualizei dados da bibi e add
# Update retention metrics for cohort rows 27, 31 and 37 on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: cohort_year 2021 (A27), period_index 4 -> num_customers 54 -> 56, retention_rate recalculated
$ws.Range("C27").Value = 56
$ws.Range("E27").Value = 0.02486678507992895

# Row 31: cohort_year 2022 (A31), period_index 3 -> num_customers 54 -> 56, retention_rate recalculated
$ws.Range("C31").Value = 56
$ws.Range("E31").Value = 0.02422145328719723

# Row 37: cohort_year 2025 (A37), period_index 0 -> num_customers and cohort_size 945 -> 948
$ws.Range("C37").Value = 948
$ws.Range("D37").Value = 948
